# Applies the crypto price/volume refresh from the
# "Updated cryptos list" GitHub Actions commit.
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in the
# sheet (leading/trailing zeros, "1.003" style grouping, "%" decorations,
# etc. all matter), so values that would otherwise be auto-coerced to a
# number on assignment are written with a leading apostrophe and the
# cell style is immediately reset to "Normal" so no numeric/quote-prefix
# formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-43: price (D) and/or volume (E) updates
$ws.Range("D2").Value = '26.996.70'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '1.827.05'
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.80%  '

$ws.Range("D5").Value = '''311.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").Value = '''1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.68%  '

$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").Value = '''0.3713'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.40%  '

$ws.Range("D9").Value = '''0.07339'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("D10").Value = '''0.8785'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.41%  '

$ws.Range("D11").Value = '''0.07890'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.47%  '

$ws.Range("D12").Value = '''19.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").Value = '1.828.69'
$ws.Range("E13").Value = '  -0.25%  '

$ws.Range("D14").Value = '''5.340'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.22%  '

$ws.Range("D15").Value = '''6.542'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.35%  '

$ws.Range("D16").Value = '''91.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.88%  '

$ws.Range("D17").Value = '''1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.49%  '

$ws.Range("D18").Value = '''0.000008833'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.22%  '

$ws.Range("D19").Value = '''1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").Value = '''14.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.52%  '

$ws.Range("D21").Value = '27.019.22'
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = '''5.095'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").Value = '''10.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.35%  '

$ws.Range("D24").Value = '2.037.69'
$ws.Range("E24").Value = '  -2.56%  '

$ws.Range("E25").Value = '  +1.09%  '

$ws.Range("D26").Value = '''1.851'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.70%  '

$ws.Range("D27").Value = '''18.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("D28").Value = '''2.049'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.46%  '

$ws.Range("D29").Value = '''5.131'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.05%  '

$ws.Range("D30").Value = '''115.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").Value = '''0.08883'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("D32").Value = '''2.959'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("D33").Value = '''0.7286'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").Value = '''4.438'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.51%  '

$ws.Range("D35").Value = '''1.131'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.65%  '

$ws.Range("D36").Value = '''2.472'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.86%  '

$ws.Range("D37").Value = '''0.01949'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.68%  '

$ws.Range("E38").Value = '  -0.62%  '

$ws.Range("D39").Value = '''0.05225'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.60%  '

$ws.Range("D40").Value = '''2.950'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.60%  '

$ws.Range("D41").Value = '''7.106'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.31%  '

$ws.Range("D42").Value = '''0.5166'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '

$ws.Range("E43").Value = '  -0.38%  '

# Rows 44-45: Decentraland and Aptos swapped order, with updated price/volume
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = '''8.175'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.59%  '

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = '''0.4842'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.24%  '

# Rows 46-51: remaining price (D) and/or volume (E) updates
$ws.Range("D46").Value = '''1.004'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = '''10.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.85%  '

$ws.Range("D48").Value = '''102.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.72%  '

$ws.Range("D49").Value = '''1.630'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("D50").Value = '''0.06199'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").Value = '''64.90'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.75%  '
